$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values for rows 13-49 (A = 11..47), which currently hold 0 or 8, to 15
for ($r = 13; $r -le 49; $r++) {
    $ws.Cells.Item($r, 2).Value = 15
}

# Extend data down to row 106.
# Rows 50-103 correspond to A = 48..101, B = 15
# Rows 104-106 correspond to A = 102..104, B = 6
for ($r = 50; $r -le 103; $r++) {
    $a = $r - 2
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = 15
}

for ($r = 104; $r -le 106; $r++) {
    $a = $r - 2
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = 6
}

# Update the view: select B2:B13 (this also clears the previous topLeftCell scroll position)
$ws.Range("B2:B13").Select()
